$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text for row 32 (SQFLite implementation note)
$ws.Range("E32").Value = "Creating SQFLite implementation for character storage + CRUD methods + routing from existing pages"

# Fill in row 33 data: start/end time, description
$ws.Range("B32").Value = 10
$ws.Range("B33").Value = 12
$ws.Range("C33").Value = 16
$ws.Range("E33").Value = "Database modifications, adding fields and modifying form"

# Update sheet view / selection to match post-edit state
$ws.Range("E33").Select()
$excel.ActiveWindow.ScrollRow = 7
